$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy header formatting (bold, border, centered) from an existing header cell
# onto the three new header cells, then set their text.
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Fill in the Wins/Losses/Ties columns for every data row (2-48).
for ($r = 2; $r -le 48; $r++) {
    $ws.Cells.Item($r, 30).Value = 100
    $ws.Cells.Item($r, 31).Value = 61
    $ws.Cells.Item($r, 32).Value = 0
}
